# "Use Generative Extraction" — reshape the invoice Report.xlsx from a
# simple 4-column table into the full generative-extraction layout:
# Invoice Number | Issue Date | Due Date | Vendor Name | Line Items |
# Total Amount | Payment Terms

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: extend the header row with the new extracted fields ---
# D1 used to read "Total"; it becomes "Vendor Name", and three new
# headers are appended after it.
$ws.Range("D1").Value = "Vendor Name"
$ws.Range("E1").Value = "Line Items"
$ws.Range("F1").Value = "Total Amount"
$ws.Range("G1").Value = "Payment Terms"

# --- Row 2: the extracted values ---
# A2 (INV-10012) is untouched.

# Issue Date / Due Date move from real dates to the extractor's raw
# string output, keeping their existing date-formatted style.
$ws.Range("B2").Value = "26/3/2021"
$ws.Range("C2").Value = "25/4/2021"

# The old single "Total" number in D2 is gone - the value now lives in
# F2 (Total Amount) instead.
$ws.Range("D2").ClearContents()

# Line Items: multi-line extracted text, wrapped.
$lineItems = "- Services, `$55.00, 10, `$550.00`n- Consulting, `$75.00, 15, `$1,125.00`n- Materials, `$123.39, 1, `$123.39"
$ws.Range("E2").Value = $lineItems
$ws.Range("E2").WrapText = $true

# Total Amount: same number as before (1699.48), but now rendered with
# a custom "R" currency format.
$ws.Range("F2").Value = 1699.48
$ws.Range("F2").NumberFormat = '"R"#,##0.00_);[Red]\("R"#,##0.00\)'

# Payment Terms: new extracted text field.
$ws.Range("G2").Value = "Please pay within 30 days using the link in your invoice email."

# The row now holds several wrapped lines of text - expand it to fit.
$ws.Rows.Item(2).RowHeight = 255
